$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "105.26"
# are not auto-converted to numbers by Excel, matching the source data
# which stores prices as literal strings.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '47.477.40'
$ws.Range('E2').Value = '  +5.62%  '
$ws.Range('D3').Value = '2.494.79'
$ws.Range('E3').Value = '  +2.96%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '323.31'
$ws.Range('E5').Value = '  +2.46%  '
$ws.Range('D6').Value = '105.26'
$ws.Range('E6').Value = '  +2.65%  '
$ws.Range('D7').Value = '0.521'
$ws.Range('E7').Value = '  +1.64%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').Value = '0.540'
$ws.Range('E9').Value = '  +3.18%  '
$ws.Range('D10').Value = '37.55'
$ws.Range('E10').Value = '  +6.18%  '
$ws.Range('D11').Value = '0.0812'
$ws.Range('E11').Value = '  +1.47%  '
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('D13').Value = '18.33'
$ws.Range('E13').Value = '  +0.52%  '
$ws.Range('D14').Value = '7.19'
$ws.Range('E14').Value = '  +3.17%  '
$ws.Range('D15').Value = '2.879.43'
$ws.Range('E15').Value = '  +2.79%  '
$ws.Range('D16').Value = '2.499.34'
$ws.Range('E16').Value = '  +2.44%  '
$ws.Range('D17').Value = '0.843'
$ws.Range('E17').Value = '  +1.10%  '
$ws.Range('D18').Value = '47.338.63'
$ws.Range('E18').Value = '  +5.52%  '
$ws.Range('D19').Value = '12.74'
$ws.Range('E19').Value = '  +4.42%  '
$ws.Range('D20').Value = '6.55'
$ws.Range('E20').Value = '  +3.08%  '
$ws.Range('D21').Value = '0.0₃0934'
$ws.Range('E21').Value = '  +1.30%  '
$ws.Range('D22').Value = '70.70'
$ws.Range('E22').Value = '  +2.98%  '
$ws.Range('D23').Value = '250.70'
$ws.Range('E23').Value = '  +3.14%  '
$ws.Range('D24').Value = '2.39'
$ws.Range('E24').Value = '  +5.86%  '
$ws.Range('D25').Value = '2.57'
$ws.Range('E25').Value = '  +3.68%  '
$ws.Range('D26').Value = '26.15'
$ws.Range('E26').Value = '  +3.96%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').Value = '10.10'
$ws.Range('E28').Value = '  +5.88%  '
$ws.Range('E29').Value = '  -2.29%  '
$ws.Range('D30').Value = '35.15'
$ws.Range('E30').Value = '  +7.49%  '
$ws.Range('D31').Value = '0.133'
$ws.Range('E31').Value = '  +7.72%  '
$ws.Range('D32').Value = '49.47'
$ws.Range('E32').Value = '  +0.83%  '
$ws.Range('D33').Value = '20.00'
$ws.Range('E33').Value = '  +0.84%  '
$ws.Range('E34').Value = '  +3.05%  '
$ws.Range('D35').Value = '0.0783'
$ws.Range('E35').Value = '  +2.95%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').Value = '4.65'
$ws.Range('E37').Value = '  +5.74%  '
$ws.Range('D38').Value = '1.94'
$ws.Range('E38').Value = '  +3.80%  '
$ws.Range('D39').Value = '2.99'
$ws.Range('E39').Value = '  +4.55%  '
$ws.Range('D41').Value = '2.23'
$ws.Range('E41').Value = '  +0.88%  '
$ws.Range('D42').Value = '120.72'
$ws.Range('E42').Value = '  -1.22%  '
$ws.Range('D43').Value = '21.44'
$ws.Range('E43').Value = '  +3.47%  '
$ws.Range('D45').Value = '1.961.26'
$ws.Range('E45').Value = '  +1.61%  '
$ws.Range('D46').Value = '2.98'
$ws.Range('E46').Value = '  +1.90%  '
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('D48').Value = '9.23'
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('D49').Value = '1.81'
$ws.Range('E49').Value = '  +1.85%  '
$ws.Range('D50').Value = '5.34'
$ws.Range('E50').Value = '  +14.45%  '
$ws.Range('D51').Value = '78.69'
$ws.Range('E51').Value = '  +3.40%  '

# Restore the default cell style on column D so no stray number format
# is left applied to the cells (matches original unstyled inline strings).
$ws.Range("D2:D51").Style = "Normal"
